{"js": "// Play implements simple max(points) strategy\n// Rename the \"Dock\" concept to \"Hand\" throughout the design doc, change\n// ExploreMove/TryExplore signatures accordingly, split the \"Dock newD = ...\"\n// line into an explicit newB.AddMove(move) call plus the renamed Hand clone\n// line, rename newD -> NewH in the TryExplore calls, and drop the stray\n// lastRenderedPageBreak field on the \"Find if a cell...\" paragraph.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// 1) \"For each of the tiles in the dock compatible ...\" -> \"...Hand...\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"For each of the tiles in the dock compatible with current playable square:\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\n    \"For each of the tiles in the Hand compatible with current playable square:\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 2) \"ExploreMove(board, dock, currentMoves, ...\" -> \"...board, Hand,...\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"(board, dock, \", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"(board, Hand, \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 3) \"Replenish user dock\" -> \"Replenish user Hand\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"Replenish user dock\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"Replenish user Hand\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 4) \"Dock is a simple set of Tile.\" -> \"Hand is a simple set of Tile.\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"Dock is a simple set of Tile.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"Hand is a simple set of Tile.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 5) \"def ExploreMove(Board b, Dock d, list<Move> ...\" signature change\n//    \"(Board b, Dock d, list\" -> \"(Board b, Hand h, list\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"(Board b, Dock d, list\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"(Board b, Hand h, list\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 6) \"Board newB = ChainedBoard(b, move)\" -> \"...ChainedBoard(b)\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"ChainedBoard(b, move)\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"ChainedBoard(b)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 7) Split \"Dock newD = d.clone \\u2013 move.tile\" into two lines:\n//      \"newB.AddMove(move)\"\n//      \"Hand newH = h.clone \\u2013 move.tile\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"Dock newD = d.clone\", { matchCase: true });\n  results.load(\"items/paragraphs\");\n  await context.sync();\n  const targetPara = results.items[0].paragraphs.getFirst();\n  targetPara.load(\"text\");\n  await context.sync();\n\n  const newPara1 = targetPara.insertParagraph(\"newB.AddMove(move)\", Word.InsertLocation.before);\n  await context.sync();\n  newPara1.insertParagraph(\"Hand newH = h.clone \\u2013 move.tile\", Word.InsertLocation.after);\n  await context.sync();\n  targetPara.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 8) \"newD\" -> \"NewH\" within the 4 TryExplore(...) calls\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"newD\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(\"NewH\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 9) \"def TryExplore (Board b, Dock d, list<Move> ...\" -> \"...Hand d,...\"\n// ---------------------------------------------------------------\n{\n  const results = body.search(\" (Board b, Dock d, list<Move> \", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\" (Board b, Hand d, list<Move> \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 10) Remove the lastRenderedPageBreak field from the \"Find if a cell is\n//     compatible with a tile\" paragraph (purely a rendering artifact).\n// ---------------------------------------------------------------\n{\n  const results = body.search(\"Find if a cell is compatible with a tile\", { matchCase: true });\n  results.load(\"items/paragraphs\");\n  await context.sync();\n  const targetPara = results.items[0].paragraphs.getFirst();\n  targetPara.load(\"text\");\n  await context.sync();\n  // Rewriting the paragraph text normalizes the run and drops any\n  // lastRenderedPageBreak field that was attached to the original run.\n  targetPara.insertText(\"Find if a cell is compatible with a tile\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Play implements simple max(points) strategy\n# Rename the \"Dock\" concept to \"Hand\" throughout the design doc, change\n# ExploreMove/TryExplore signatures accordingly, split the \"Dock newD = ...\"\n# line into an explicit newB.AddMove(move) call plus the renamed Hand clone\n# line, rename newD -> NewH in the TryExplore calls, and drop the stray\n# lastRenderedPageBreak field on the \"Find if a cell...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------\n# 1) \"For each of the tiles in the dock compatible ...\" -> \"...Hand...\"\n#    (paragraph 14, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(14).Range\n$r.Find.Execute(\"the dock compatible\", $true, $false, $false, $false, $false, $true, 1, $false, \"the Hand compatible\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 2) \"ExploreMove(board, dock, currentMoves, ...\" -> \"...board, Hand,...\"\n#    (paragraph 16, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(16).Range\n$r.Find.Execute(\"(board, dock, \", $true, $false, $false, $false, $false, $true, 1, $false, \"(board, Hand, \", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 3) \"Replenish user dock\" -> \"Replenish user Hand\"\n#    (paragraph 20, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(20).Range\n$r.Find.Execute(\"Replenish user dock\", $true, $false, $false, $false, $false, $true, 1, $false, \"Replenish user Hand\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 4) \"Dock is a simple set of Tile.\" -> \"Hand is a simple set of Tile.\"\n#    (paragraph 23, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(23).Range\n$r.Find.Execute(\"Dock is a simple set of Tile.\", $true, $false, $false, $false, $false, $true, 1, $false, \"Hand is a simple set of Tile.\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 5) \"def ExploreMove(Board b, Dock d, list<Move> ...\" signature change\n#    \"Dock d, list\" -> \"Hand h, list\"  (paragraph 27, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(27).Range\n$r.Find.Execute(\"Dock d, list\", $true, $false, $false, $false, $false, $true, 1, $false, \"Hand h, list\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 6) \"Board newB = ChainedBoard(b, move)\" -> \"...ChainedBoard(b)\"\n#    (paragraph 28, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(28).Range\n$r.Find.Execute(\"ChainedBoard(b, move)\", $true, $false, $false, $false, $false, $true, 1, $false, \"ChainedBoard(b)\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 8) \"newD\" -> \"NewH\" within the 4 TryExplore(...) calls\n#    (paragraphs 33, 34, 36, 37 -- 1-based; these come *after* the\n#    split point so handle them before inserting the new paragraph)\n# ---------------------------------------------------------------\nforeach ($idx in 33, 34, 36, 37) {\n  $r = $d.Paragraphs.Item($idx).Range\n  $r.Find.Execute(\"newD\", $true, $false, $false, $false, $false, $true, 1, $false, \"NewH\", 2) | Out-Null\n}\n\n# ---------------------------------------------------------------\n# 9) \"def TryExplore (Board b, Dock d, list<Move> ...\" -> \"...Hand d,...\"\n#    (paragraph 39, 1-based)\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(39).Range\n$r.Find.Execute(\"Dock d, list\", $true, $false, $false, $false, $false, $true, 1, $false, \"Hand d, list\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 10) Remove the lastRenderedPageBreak field from the \"Find if a cell is\n#     compatible with a tile\" paragraph (paragraph 52, 1-based) -- simply\n#     re-stamping the text with itself rebuilds the run without the field.\n# ---------------------------------------------------------------\n$r = $d.Paragraphs.Item(52).Range\n$r.Find.Execute(\"Find if a cell is compatible with a tile\", $true, $false, $false, $false, $false, $true, 1, $false, \"Find if a cell is compatible with a tile\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 7) Split \"Dock newD = d.clone - move.tile\" (paragraph 29, 1-based) into:\n#      \"newB.AddMove(move)\"\n#      \"Hand newH = h.clone - move.tile\"\n#    Do this last since it changes paragraph count / shifts indices.\n# ---------------------------------------------------------------\n$p28 = $d.Paragraphs.Item(28)    # \"Board newB = ChainedBoard(b)\"\n$p28.Range.InsertParagraphAfter() | Out-Null\n$newPara = $d.Paragraphs.Item(29)\n$newPara.Range.InsertAfter(\"newB.AddMove(move)\")\n\n$origPara = $d.Paragraphs.Item(30)\n$r1 = $origPara.Range\n$r1.Find.Execute(\"Dock newD\", $true, $false, $false, $false, $false, $true, 1, $false, \"Hand newH\", 2) | Out-Null\n$r2 = $d.Paragraphs.Item(30).Range\n$r2.Find.Execute(\"d.clone\", $true, $false, $false, $false, $false, $true, 1, $false, \"h.clone\", 2) | Out-Null\n"}
